# Fix prisma.fi / tokmanni EAN control logic:
# Row 4 (EAN 3606480504662) is a duplicate/erroneous entry that needs to be
# removed entirely; all rows below it shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 4, shifting the rows below it upward.
$ws.Rows("4:4").Delete()

# Update the active selection to match the post-edit workbook state.
$ws.Range("D10").Select()
